# Refresh the cryptos price/volume table (GitHub Actions style data update).
# Most D-column price strings are plain text in the source sheet (e.g. "42.975.72",
# "0.998"); values with a leading apostrophe are assigned so Excel keeps them as text
# instead of silently parsing them into numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.006.09'
$ws.Range('E2').Value = '  +0.00%  '

$ws.Range('D3').Value = '2.561.04'
$ws.Range('E3').Value = '  +0.82%  '

$ws.Range('D4').Value = "'1.00"
$ws.Range('E4').Value = '  -0.05%  '

$ws.Range('D5').Value = "'313.69"
$ws.Range('E5').Value = '  -1.28%  '

$ws.Range('D6').Value = "'96.59"
$ws.Range('E6').Value = '  +0.92%  '

$ws.Range('E7').Value = '  -0.38%  '

$ws.Range('D8').Value = "'1.00"
$ws.Range('E8').Value = '  -0.14%  '

$ws.Range('D9').Value = "'0.540"
$ws.Range('E9').Value = '  +1.35%  '

$ws.Range('E10').Value = '  -2.31%  '

$ws.Range('D11').Value = "'0.0812"
$ws.Range('E11').Value = '  -0.12%  '

$ws.Range('E12').Value = '  -2.65%  '

$ws.Range('B13').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C13').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D13').Value = '2.955.09'
$ws.Range('E13').Value = '  +0.61%  '

$ws.Range('B14').Value = 'TRON'
$ws.Range('C14').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D14').Value = "'0.108"
$ws.Range('E14').Value = '  -4.71%  '

$ws.Range('D15').Value = '2.555.89'
$ws.Range('E15').Value = '  -0.83%  '

$ws.Range('D16').Value = "'15.05"
$ws.Range('E16').Value = '  -2.87%  '

$ws.Range('E17').Value = '  -1.29%  '

$ws.Range('D18').Value = '42.999.73'
$ws.Range('E18').Value = '  -0.05%  '

$ws.Range('E19').Value = '  +2.49%  '

$ws.Range('D20').Value = "'12.52"
$ws.Range('E20').Value = '  -4.22%  '

$ws.Range('D21').Value = '0.0₃0958'
$ws.Range('E21').Value = '  -1.30%  '

$ws.Range('D22').Value = "'69.10"
$ws.Range('E22').Value = '  -1.87%  '

$ws.Range('D23').Value = "'252.09"
$ws.Range('E23').Value = '  -0.14%  '

$ws.Range('D24').Value = "'2.94"
$ws.Range('E24').Value = '  -1.02%  '

$ws.Range('D25').Value = "'2.07"
$ws.Range('E25').Value = '  +2.24%  '

$ws.Range('D26').Value = "'26.72"
$ws.Range('E26').Value = '  -1.36%  '

$ws.Range('E27').Value = '  -0.07%  '

$ws.Range('E28').Value = '  +0.14%  '

$ws.Range('D29').Value = "'39.96"
$ws.Range('E29').Value = '  +0.04%  '

$ws.Range('E30').Value = '  -0.34%  '

$ws.Range('D31').Value = "'5.80"
$ws.Range('E31').Value = '  -4.77%  '

$ws.Range('D32').Value = "'154.33"
$ws.Range('E32').Value = '  -0.65%  '

$ws.Range('E33').Value = '  +1.81%  '

$ws.Range('E34').Value = '  +1.66%  '

$ws.Range('B35').Value = 'ARBITRUM'
$ws.Range('C35').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D35').Value = "'2.12"
$ws.Range('E35').Value = '  -0.86%  '

$ws.Range('B36').Value = 'WEMIXToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D36').Value = "'2.70"
$ws.Range('E36').Value = '  +2.66%  '

$ws.Range('D37').Value = "'19.02"
$ws.Range('E37').Value = '  +0.33%  '

$ws.Range('E38').Value = '  -1.38%  '

$ws.Range('E39').Value = '  +6.48%  '

$ws.Range('E40').Value = '  -0.85%  '

$ws.Range('D41').Value = "'22.56"
$ws.Range('E41').Value = '  -7.57%  '

$ws.Range('D42').Value = "'3.93"
$ws.Range('E42').Value = '  +3.02%  '

$ws.Range('E43').Value = '  +0.14%  '

$ws.Range('E44').Value = '  -0.15%  '

$ws.Range('E45').Value = '  -3.82%  '

$ws.Range('D46').Value = '2.002.64'
$ws.Range('E46').Value = '  -0.88%  '

$ws.Range('E47').Value = '  +1.42%  '

$ws.Range('D48').Value = "'83.08"
$ws.Range('E48').Value = '  -3.25%  '

$ws.Range('D49').Value = '2.804.45'
$ws.Range('E49').Value = '  +0.43%  '

$ws.Range('D50').Value = "'74.13"
$ws.Range('E50').Value = '  -0.24%  '

$ws.Range('B51').Value = 'Algorand'
$ws.Range('C51').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D51').Value = "'0.193"
$ws.Range('E51').Value = '  +1.60%  '
